$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Move the "Total" row from 36 down to 39 (copy formats first, then
#        content), before we overwrite row 36 with new timesheet data.
$ws.Range("B36:C36").Copy()
$ws.Range("B39").PasteSpecial(-4122)
$ws.Range("B39").Value = "Total"
$ws.Range("C39").Formula = "=SUM(C5:C38)"

# --- 2) Update the wording of the two activity descriptions (shared by
#        several rows). Updating every cell that references the old text to
#        the same new text lets the engine fold them back onto a single
#        shared-string slot instead of spawning a duplicate.
$ws.Range("B32").Value = "Début d'un gros refactoring de l'entierté du code + Javadoc"
$ws.Range("B33").Value = "Suite refactoring et débugging et Javadoc"
$ws.Range("B34").Value = "Suite refactoring et débugging et Javadoc"
$ws.Range("B35").Value = "Suite refactoring et débugging et Javadoc"

# --- 3) Row 35: hours value changes from 1.5 to 10.
$ws.Range("C35").Value = 10
$ws.Range("C34").Copy()
$ws.Range("C35").PasteSpecial(-4122)

# --- 4) New row 36: 19 May 2018, same activity, 1 hour.
$ws.Range("A35:C35").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A36").Value = 43239
$ws.Range("B36").Value = "Suite refactoring et débugging et Javadoc"
$ws.Range("C36").Value = 1

# --- 5) New row 37: 20 May 2018, same activity, hours left blank.
$ws.Range("A35:C35").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A37").Value = 43240
$ws.Range("B37").Value = "Suite refactoring et débugging et Javadoc"

# --- 6) New row 38: blank spacer row (keeps the date/activity formatting,
#        but the hours cell uses the centered-without-border-dup style).
$ws.Range("A35:B35").Copy()
$ws.Range("A38").PasteSpecial(-4122)
$ws.Range("C23").Copy()
$ws.Range("C38").PasteSpecial(-4122)

# --- 7) Selection, matching what Excel leaves behind after this edit.
$ws.Range("E38").Select()
